# Update the invoice row with the latest test data and restore the
# page setup (paper size / orientation) used when testing the Excel
# activities.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New invoice number
$ws.Range("B2").Value = 8669

# New customer name
$ws.Range("C2").Value = "Lance"

# Date typed in as plain dd/mm/yyyy text rather than a real date value
$ws.Range("D2").Style = "Standard"
$ws.Range("D2").Value = "22/12/2024"

# Page setup used for printing/testing
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
